$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 131144039
$ws.Range("B3").Value = 56762
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 100092
$ws.Range("F3").Value = "Större brunfladdermus"
$ws.Range("G3").Value = "Nyctalus noctula"
$ws.Range("H3").Value = "(Schreber, 1774)"
$ws.Range("I3").Value = "'3"
$ws.Cells.Item(3,9).Style = "Normal"
$ws.Range("J3").Value = "registreringar"
$ws.Range("M3").Value = "födosökande"
$ws.Range("N3").Value = "autobox med tidsexpansion"
$ws.Range("P3").Value = "Södra Atriumvägen, Upl"
$ws.Range("Q3").Value = 686956
$ws.Range("R3").Value = 6601143
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Stockholm"
$ws.Range("U3").Value = "Österåker"
$ws.Range("V3").Value = "Uppland"
$ws.Range("W3").Value = "Österåker"
$ws.Range("Y3").Value = "'2025-08-21"
$ws.Cells.Item(3,25).Style = "Normal"
$ws.Range("Z3").Value = "22:00"
$ws.Range("AA3").Value = "'2025-08-23"
$ws.Cells.Item(3,27).Style = "Normal"
$ws.Range("AB3").Value = "06:00"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = "Martin Berg"
$ws.Range("AX3").Value = "Martin Berg"

# Row 4
$ws.Range("A4").Value = 131143868
$ws.Range("B4").Value = 56762
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 100092
$ws.Range("F4").Value = "Större brunfladdermus"
$ws.Range("G4").Value = "Nyctalus noctula"
$ws.Range("H4").Value = "(Schreber, 1774)"
$ws.Range("I4").Value = "'12"
$ws.Cells.Item(4,9).Style = "Normal"
$ws.Range("J4").Value = "registreringar"
$ws.Range("M4").Value = "födosökande"
$ws.Range("N4").Value = "autobox med tidsexpansion"
$ws.Range("P4").Value = "Båtstorps förskola, Upl"
$ws.Range("Q4").Value = 686690
$ws.Range("R4").Value = 6601110
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Stockholm"
$ws.Range("U4").Value = "Österåker"
$ws.Range("V4").Value = "Uppland"
$ws.Range("W4").Value = "Österåker"
$ws.Range("Y4").Value = "'2025-07-23"
$ws.Cells.Item(4,25).Style = "Normal"
$ws.Range("Z4").Value = "22:00"
$ws.Range("AA4").Value = "'2025-07-25"
$ws.Cells.Item(4,27).Style = "Normal"
$ws.Range("AB4").Value = "06:00"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "Martin Berg"
$ws.Range("AX4").Value = "Martin Berg"

# Row 5
$ws.Range("A5").Value = 131144044
$ws.Range("B5").Value = 56748
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 205998
$ws.Range("F5").Value = "Nordfladdermus"
$ws.Range("G5").Value = "Eptesicus nilssonii"
$ws.Range("H5").Value = "(A.Keyserling & Blasius, 1839)"
$ws.Range("I5").Value = "'1"
$ws.Cells.Item(5,9).Style = "Normal"
$ws.Range("J5").Value = "registreringar"
$ws.Range("M5").Value = "födosökande"
$ws.Range("N5").Value = "autobox med tidsexpansion"
$ws.Range("P5").Value = "Södra Atriumvägen, Upl"
$ws.Range("Q5").Value = 686956
$ws.Range("R5").Value = 6601143
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Stockholm"
$ws.Range("U5").Value = "Österåker"
$ws.Range("V5").Value = "Uppland"
$ws.Range("W5").Value = "Österåker"
$ws.Range("Y5").Value = "'2025-07-23"
$ws.Cells.Item(5,25).Style = "Normal"
$ws.Range("Z5").Value = "22:00"
$ws.Range("AA5").Value = "'2025-07-25"
$ws.Cells.Item(5,27).Style = "Normal"
$ws.Range("AB5").Value = "06:00"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = "Martin Berg"
$ws.Range("AX5").Value = "Martin Berg"

# Row 6
$ws.Range("A6").Value = 131143867
$ws.Range("B6").Value = 56748
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 205998
$ws.Range("F6").Value = "Nordfladdermus"
$ws.Range("G6").Value = "Eptesicus nilssonii"
$ws.Range("H6").Value = "(A.Keyserling & Blasius, 1839)"
$ws.Range("I6").Value = "'15"
$ws.Cells.Item(6,9).Style = "Normal"
$ws.Range("J6").Value = "registreringar"
$ws.Range("M6").Value = "födosökande"
$ws.Range("N6").Value = "autobox med tidsexpansion"
$ws.Range("P6").Value = "Båtstorps förskola, Upl"
$ws.Range("Q6").Value = 686690
$ws.Range("R6").Value = 6601110
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Stockholm"
$ws.Range("U6").Value = "Österåker"
$ws.Range("V6").Value = "Uppland"
$ws.Range("W6").Value = "Österåker"
$ws.Range("Y6").Value = "'2025-07-23"
$ws.Cells.Item(6,25).Style = "Normal"
$ws.Range("Z6").Value = "22:00"
$ws.Range("AA6").Value = "'2025-07-25"
$ws.Cells.Item(6,27).Style = "Normal"
$ws.Range("AB6").Value = "06:00"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "Martin Berg"
$ws.Range("AX6").Value = "Martin Berg"

# Row 7
$ws.Range("A7").Value = 131143865
$ws.Range("B7").Value = 56769
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 206002
$ws.Range("F7").Value = "Brunlångöra"
$ws.Range("G7").Value = "Plecotus auritus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("I7").Value = "'1"
$ws.Cells.Item(7,9).Style = "Normal"
$ws.Range("J7").Value = "registreringar"
$ws.Range("M7").Value = "födosökande"
$ws.Range("N7").Value = "autobox med tidsexpansion"
$ws.Range("P7").Value = "Båtstorps förskola, Upl"
$ws.Range("Q7").Value = 686690
$ws.Range("R7").Value = 6601110
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = "Stockholm"
$ws.Range("U7").Value = "Österåker"
$ws.Range("V7").Value = "Uppland"
$ws.Range("W7").Value = "Österåker"
$ws.Range("Y7").Value = "'2025-07-23"
$ws.Cells.Item(7,25).Style = "Normal"
$ws.Range("Z7").Value = "22:00"
$ws.Range("AA7").Value = "'2025-07-25"
$ws.Cells.Item(7,27).Style = "Normal"
$ws.Range("AB7").Value = "06:00"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = "Martin Berg"
$ws.Range("AX7").Value = "Martin Berg"

# Row 8
$ws.Range("A8").Value = 131144032
$ws.Range("B8").Value = 56767
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 205995
$ws.Range("F8").Value = "Dvärgpipistrell"
$ws.Range("G8").Value = "Pipistrellus pygmaeus"
$ws.Range("H8").Value = "(W.E.Leach, 1825)"
$ws.Range("I8").Value = "'2"
$ws.Cells.Item(8,9).Style = "Normal"
$ws.Range("J8").Value = "registreringar"
$ws.Range("M8").Value = "födosökande"
$ws.Range("N8").Value = "autobox med tidsexpansion"
$ws.Range("P8").Value = "Södra Atriumvägen, Upl"
$ws.Range("Q8").Value = 686956
$ws.Range("R8").Value = 6601143
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = "Stockholm"
$ws.Range("U8").Value = "Österåker"
$ws.Range("V8").Value = "Uppland"
$ws.Range("W8").Value = "Österåker"
$ws.Range("Y8").Value = "'2025-08-21"
$ws.Cells.Item(8,25).Style = "Normal"
$ws.Range("Z8").Value = "22:00"
$ws.Range("AA8").Value = "'2025-08-23"
$ws.Cells.Item(8,27).Style = "Normal"
$ws.Range("AB8").Value = "06:00"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = "Martin Berg"
$ws.Range("AX8").Value = "Martin Berg"
